# Apply "zero_before_threshold" recalculation to the Step3_DataPts_* sheets.
#
# Zeroing signal values before the noise threshold / first rise point shifts
# the detected First_Noticeable_Increase_Index (col C) earlier for several
# segments, which in turn changes the First_Noticeable_Increase_Cumulative_Value
# (col E) and the derived Pulse_Width (col G = Point_Exceeds_Index - First_Noticeable_Increase_Index).
# Point_Exceeds_Index (col D) / Point_Exceeds_Cumulative_Value (col F) are unaffected.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Step3_DataPts_0.5",
    "Step3_DataPts_0.7",
    "Step3_DataPts_0.8",
    "Step3_DataPts_0.9"
)

# New First_Noticeable_Increase_Index (C) and First_Noticeable_Increase_Cumulative_Value (E)
# per data row (segment). These are identical across the four threshold sheets because the
# underlying segment signal does not depend on the Intensity_Threshold column.
$rowUpdates = @{
    2 = @{ C = 87; E = 0.0009445236536463978 }
    3 = @{ C = 87; E = 0.002291213658635639 }
    4 = @{ C = 88; E = 0.001525559387765866 }
    5 = @{ C = 87; E = 0.001613222062089176 }
    6 = @{ C = 87; E = 0.001729082465001621 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $rowUpdates.Keys) {
        $update = $rowUpdates[$row]

        # Point_Exceeds_Index (column D) stays as-is; Pulse_Width (column G)
        # is recomputed from the new First_Noticeable_Increase_Index.
        $pointExceedsIndex = $ws.Cells.Item($row, 4).Value2

        $ws.Cells.Item($row, 3).Value = $update.C
        $ws.Cells.Item($row, 5).Value = $update.E
        $ws.Cells.Item($row, 7).Value = $pointExceedsIndex - $update.C
    }
}
